$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D73: 93.3 -> 98.3
$ws.Range("D73").Value = 98.3

# Swap rows 80 and 81 (CHIR <-> DAPT blocks), including styles
$c80 = $ws.Range("C80:F80").Value
$c81 = $ws.Range("C81:F81").Value
$s80 = $ws.Range("C80:F80").Style
$s81 = $ws.Range("C81:F81").Style

$ws.Range("C80:F80").Value = $c81
$ws.Range("C80:F80").Style = $s81
$ws.Range("C81:F81").Value = $c80
$ws.Range("C81:F81").Style = $s80

# Swap rows 95 and 96
$c95 = $ws.Range("C95:F95").Value
$c96 = $ws.Range("C96:F96").Value
$s95 = $ws.Range("C95:F95").Style
$s96 = $ws.Range("C96:F96").Style

$ws.Range("C95:F95").Value = $c96
$ws.Range("C95:F95").Style = $s96
$ws.Range("C96:F96").Value = $c95
$ws.Range("C96:F96").Style = $s95

# Rows 88-89 value + style
$ws.Range("D88").Value = 98.3
$ws.Range("C88:D89").Style = $ws.Range("C109:D110").Style

# Row 100-101
$ws.Range("D100").Value = 98.3

# Update styles for rows 88,89,100,101,109,110 C/D columns to use the "new" style (font default instead of bold red)
